$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.091.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.768.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5226"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2758"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.57"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06195"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.775.63"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07025"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.78"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6462"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.530"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.14"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "26.126.26"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.68"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006763"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.000.15"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.076"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.467"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.190"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.93"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.496"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.855"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.21"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.14"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08418"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.712"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.461"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04474"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.657"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.008"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6055"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.752"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01590"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.990"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.69%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.86"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3883"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7524"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.947"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05519"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.378"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1119"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.67"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.79%  "
